$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 365.30768
$ws.Range("I2").Value = 151.5
$ws.Range("K2").Value = 151.5
$ws.Range("M2").Value = -38.5

$ws.Range("H64").Value = 76927090
$ws.Range("I64").Value = 333335680
$ws.Range("J64").Value = 4513
$ws.Range("K64").Value = 333335680
$ws.Range("L64").Value = 4513
$ws.Range("M64").Value = -333335432
$ws.Range("N64").Value = -5009

$ws.Range("H67").Value = 76927090
$ws.Range("I67").Value = 333335680
$ws.Range("J67").Value = 4513
$ws.Range("K67").Value = 333335680
$ws.Range("L67").Value = 4513
$ws.Range("M67").Value = -333334822
$ws.Range("N67").Value = -6229.7334

$ws.Range("H69").Value = 4020.7292
$ws.Range("J69").Value = 3993.3872
$ws.Range("L69").Value = 11980.1616
$ws.Range("N69").Value = -13728.1616

$ws.Range("H72").Value = 4020.7292
$ws.Range("J72").Value = 3993.3872
$ws.Range("L72").Value = 35940.4848
$ws.Range("N72").Value = -44676.4848

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H76").Value = 5005.4443
$ws.Range("I76").Value = 3900
$ws.Range("J76").Value = 5889.8
$ws.Range("K76").Value = 3900
$ws.Range("L76").Value = 5889.8
$ws.Range("M76").Value = -3585
$ws.Range("N76").Value = -6519.8

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H79").Value = 5005.4443
$ws.Range("I79").Value = 3900
$ws.Range("J79").Value = 5889.8
$ws.Range("K79").Value = 3900
$ws.Range("L79").Value = 5889.8
$ws.Range("M79").Value = -2808
$ws.Range("N79").Value = -8073.8

$ws.Range("H80").Value = 696.4643
$ws.Range("I80").Value = 351.2857
$ws.Range("J80").Value = 1041.6428
$ws.Range("K80").Value = 1053.8571
$ws.Range("L80").Value = 3124.9284
$ws.Range("M80").Value = -55.85710000000017
$ws.Range("N80").Value = -5120.928400000001

$ws.Range("H83").Value = 696.4643
$ws.Range("I83").Value = 351.2857
$ws.Range("J83").Value = 1041.6428
$ws.Range("K83").Value = 3161.5713
$ws.Range("L83").Value = 9374.7852
$ws.Range("M83").Value = 1830.4287
$ws.Range("N83").Value = -19358.7852

$ws.Range("H92").Value = 454.44116
$ws.Range("I92").Value = 236.12
$ws.Range("J92").Value = 1060.8889
$ws.Range("K92").Value = 236.12
$ws.Range("L92").Value = 1060.8889
$ws.Range("M92").Value = 1011.88
$ws.Range("N92").Value = -3556.8889

$ws.Range("H137").Value = 3840.9285
$ws.Range("I137").Value = 882.4091
$ws.Range("J137").Value = 7095.3
$ws.Range("K137").Value = 2647.2273
$ws.Range("L137").Value = 21285.9
$ws.Range("M137").Value = -97.22730000000001
$ws.Range("N137").Value = -26385.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3807.6365
$ws.Range("I74").Value = 506
$ws.Range("J74").Value = 14125.25
$ws.Range("K74").Value = 506
$ws.Range("L74").Value = 14125.25
$ws.Range("M74").Value = 368
$ws.Range("N74").Value = -15873.25

$ws.Range("H77").Value = 3807.6365
$ws.Range("I77").Value = 506
$ws.Range("J77").Value = 14125.25
$ws.Range("K77").Value = 2530
$ws.Range("L77").Value = 70626.25
$ws.Range("M77").Value = 1838
$ws.Range("N77").Value = -79362.25

$ws.Range("H97").Value = 2590
$ws.Range("I97").Value = 710.4737
$ws.Range("J97").Value = 6161.1
$ws.Range("K97").Value = 710.4737
$ws.Range("L97").Value = 6161.1
$ws.Range("M97").Value = -214.4737
$ws.Range("N97").Value = -7153.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 772492.56
$ws.Range("I86").Value = 1900.8334
$ws.Range("J86").Value = 1432999.8
$ws.Range("K86").Value = 1900.8334
$ws.Range("L86").Value = 1432999.8
$ws.Range("M86").Value = -777.8334
$ws.Range("N86").Value = -1435245.8

$ws.Range("H89").Value = 772492.56
$ws.Range("I89").Value = 1900.8334
$ws.Range("J89").Value = 1432999.8
$ws.Range("K89").Value = 9504.167
$ws.Range("L89").Value = 7164999
$ws.Range("M89").Value = -3888.166999999999
$ws.Range("N89").Value = -7176231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32296386
$ws.Range("I31").Value = 71431070
$ws.Range("J31").Value = 67820.94
$ws.Range("K31").Value = 71431070
$ws.Range("L31").Value = 67820.94
$ws.Range("M31").Value = -71430775
$ws.Range("N31").Value = -68410.94

$ws.Range("H34").Value = 32296386
$ws.Range("I34").Value = 71431070
$ws.Range("J34").Value = 67820.94
$ws.Range("K34").Value = 71431070
$ws.Range("L34").Value = 67820.94
$ws.Range("M34").Value = -71430868
$ws.Range("N34").Value = -68224.94

$ws.Range("H62").Value = 4046.1282
$ws.Range("J62").Value = 2539.8
$ws.Range("L62").Value = 2539.8
$ws.Range("N62").Value = -3787.8

$ws.Range("H65").Value = 4046.1282
$ws.Range("J65").Value = 2539.8
$ws.Range("L65").Value = 12699
$ws.Range("N65").Value = -18939

$ws.Range("H107").Value = 510.04
$ws.Range("I107").Value = 507.3158
$ws.Range("J107").Value = 518.6667
$ws.Range("K107").Value = 507.3158
$ws.Range("L107").Value = 518.6667
$ws.Range("M107").Value = 1412.6842
$ws.Range("N107").Value = -4358.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 57.333332
$ws.Range("I23").Value = 56.833332
$ws.Range("J23").Value = 57.833332
$ws.Range("K23").Value = 170.499996
$ws.Range("L23").Value = 173.499996
$ws.Range("M23").Value = 64.50000399999999
$ws.Range("N23").Value = -643.499996

$ws.Range("H132").Value = 1143.1111
$ws.Range("J132").Value = 1325
$ws.Range("L132").Value = 11925
$ws.Range("N132").Value = -16985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4958.4546
$ws.Range("I70").Value = 4803.2
$ws.Range("J70").Value = 5087.8335
$ws.Range("K70").Value = 4803.2
$ws.Range("L70").Value = 5087.8335
$ws.Range("M70").Value = -4533.2
$ws.Range("N70").Value = -5627.8335

$ws.Range("H73").Value = 4958.4546
$ws.Range("I73").Value = 4803.2
$ws.Range("J73").Value = 5087.8335
$ws.Range("K73").Value = 4803.2
$ws.Range("L73").Value = 5087.8335
$ws.Range("M73").Value = -3867.2
$ws.Range("N73").Value = -6959.8335

$ws.Range("H97").Value = 1249.5834
$ws.Range("I97").Value = 1134.5
$ws.Range("J97").Value = 1825
$ws.Range("K97").Value = 1134.5
$ws.Range("L97").Value = 1825
$ws.Range("M97").Value = -638.5
$ws.Range("N97").Value = -2817

$ws.Range("H122").Value = 2623.8333
$ws.Range("I122").Value = 2685.875
$ws.Range("J122").Value = 2499.75
$ws.Range("K122").Value = 8057.625
$ws.Range("L122").Value = 7499.25
$ws.Range("M122").Value = -5607.625
$ws.Range("N122").Value = -12399.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 28500
$ws.Range("I74").Value = 14000
$ws.Range("K74").Value = 14000
$ws.Range("M74").Value = -13002

$ws.Range("H75").Value = 21000
$ws.Range("I75").Value = 14000
$ws.Range("J75").Value = 28000
$ws.Range("K75").Value = 14000
$ws.Range("L75").Value = 28000
$ws.Range("M75").Value = -13064
$ws.Range("N75").Value = -29872

$ws.Range("H77").Value = 28500
$ws.Range("I77").Value = 14000
$ws.Range("K77").Value = 42000
$ws.Range("M77").Value = -37008

$ws.Range("H78").Value = 21000
$ws.Range("I78").Value = 14000
$ws.Range("J78").Value = 28000
$ws.Range("K78").Value = 42000
$ws.Range("L78").Value = 84000
$ws.Range("M78").Value = -37320
$ws.Range("N78").Value = -93360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H100").Value = 356.42105
$ws.Range("I100").Value = 329.625
$ws.Range("J100").Value = 499.33334
$ws.Range("K100").Value = 659.25
$ws.Range("L100").Value = 998.66668
$ws.Range("M100").Value = -118.25
$ws.Range("N100").Value = -2080.66668

